# Apply the Natmi re-run (Dr Hou advice) edit:
# - Adds a third sending/target cluster "ECs" (shared string)
# - Expands the data block from 4 rows (2x2 clusters) to 9 rows (3x3 clusters)
# - Updates all numeric columns (E..T) accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Vtn"
$ws.Cells.Item(2,3).Value2 = "Itgb8"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 2
$ws.Cells.Item(2,6).Value2 = 0.6666666666666666
$ws.Cells.Item(2,7).Value2 = 7.134618
$ws.Cells.Item(2,8).Value2 = 21.403854
$ws.Cells.Item(2,9).Value2 = 0.0965317920926077
$ws.Cells.Item(2,10).Value2 = 0.0965317920926077
$ws.Cells.Item(2,11).Value2 = 1
$ws.Cells.Item(2,12).Value2 = 0.3333333333333333
$ws.Cells.Item(2,13).Value2 = 0.1126243333333333
$ws.Cells.Item(2,14).Value2 = 0.337873
$ws.Cells.Item(2,15).Value2 = 0.01082936903163217
$ws.Cells.Item(2,16).Value2 = 0.01082936903163217
$ws.Cells.Item(2,17).Value2 = 0.8035315958379999
$ws.Cells.Item(2,18).Value2 = 7.231784362541999
$ws.Cells.Item(2,19).Value2 = 0.001045378399855641
$ws.Cells.Item(2,20).Value2 = 0.001045378399855641

# Row 3
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Vtn"
$ws.Cells.Item(3,3).Value2 = "Itgb8"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 2
$ws.Cells.Item(3,6).Value2 = 0.6666666666666666
$ws.Cells.Item(3,7).Value2 = 7.134618
$ws.Cells.Item(3,8).Value2 = 21.403854
$ws.Cells.Item(3,9).Value2 = 0.0965317920926077
$ws.Cells.Item(3,10).Value2 = 0.0965317920926077
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 3.323421
$ws.Cells.Item(3,14).Value2 = 9.970263
$ws.Cells.Item(3,15).Value2 = 0.3195628457125252
$ws.Cells.Item(3,16).Value2 = 0.3195628457125252
$ws.Cells.Item(3,17).Value2 = 23.711339288178
$ws.Cells.Item(3,18).Value2 = 213.402053593602
$ws.Cells.Item(3,19).Value2 = 0.03084797418284356
$ws.Cells.Item(3,20).Value2 = 0.03084797418284356

# Row 4
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Vtn"
$ws.Cells.Item(4,3).Value2 = "Itgb8"
$ws.Cells.Item(4,4).Value2 = "sCs"
$ws.Cells.Item(4,5).Value2 = 2
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,7).Value2 = 7.134618
$ws.Cells.Item(4,8).Value2 = 21.403854
$ws.Cells.Item(4,9).Value2 = 0.0965317920926077
$ws.Cells.Item(4,10).Value2 = 0.0965317920926077
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 6.963852666666667
$ws.Cells.Item(4,14).Value2 = 20.891558
$ws.Cells.Item(4,15).Value2 = 0.6696077852558425
$ws.Cells.Item(4,16).Value2 = 0.6696077852558425
$ws.Cells.Item(4,17).Value2 = 49.684428584948
$ws.Cells.Item(4,18).Value2 = 447.159857264532
$ws.Cells.Item(4,19).Value2 = 0.0646384395099085
$ws.Cells.Item(4,20).Value2 = 0.0646384395099085

# Row 5
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Vtn"
$ws.Cells.Item(5,3).Value2 = "Itgb8"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 17.50798033333334
$ws.Cells.Item(5,8).Value2 = 52.52394100000001
$ws.Cells.Item(5,9).Value2 = 0.2368839813846793
$ws.Cells.Item(5,10).Value2 = 0.2368839813846794
$ws.Cells.Item(5,11).Value2 = 1
$ws.Cells.Item(5,12).Value2 = 0.3333333333333333
$ws.Cells.Item(5,13).Value2 = 0.1126243333333333
$ws.Cells.Item(5,14).Value2 = 0.337873
$ws.Cells.Item(5,15).Value2 = 0.01082936903163217
$ws.Cells.Item(5,16).Value2 = 0.01082936903163217
$ws.Cells.Item(5,17).Value2 = 1.971824613054778
$ws.Cells.Item(5,18).Value2 = 17.746421517493
$ws.Cells.Item(5,19).Value2 = 0.002565304052096977
$ws.Cells.Item(5,20).Value2 = 0.002565304052096977

# Row 6
$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(6,2).Value2 = "Vtn"
$ws.Cells.Item(6,3).Value2 = "Itgb8"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 17.50798033333334
$ws.Cells.Item(6,8).Value2 = 52.52394100000001
$ws.Cells.Item(6,9).Value2 = 0.2368839813846793
$ws.Cells.Item(6,10).Value2 = 0.2368839813846794
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 3.323421
$ws.Cells.Item(6,14).Value2 = 9.970263
$ws.Cells.Item(6,15).Value2 = 0.3195628457125252
$ws.Cells.Item(6,16).Value2 = 0.3195628457125252
$ws.Cells.Item(6,17).Value2 = 58.18638950738701
$ws.Cells.Item(6,18).Value2 = 523.6775055664831
$ws.Cells.Item(6,19).Value2 = 0.07569931919500099
$ws.Cells.Item(6,20).Value2 = 0.07569931919500099

# Row 7
$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Vtn"
$ws.Cells.Item(7,3).Value2 = "Itgb8"
$ws.Cells.Item(7,4).Value2 = "sCs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 17.50798033333334
$ws.Cells.Item(7,8).Value2 = 52.52394100000001
$ws.Cells.Item(7,9).Value2 = 0.2368839813846793
$ws.Cells.Item(7,10).Value2 = 0.2368839813846794
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 6.963852666666667
$ws.Cells.Item(7,14).Value2 = 20.891558
$ws.Cells.Item(7,15).Value2 = 0.6696077852558425
$ws.Cells.Item(7,16).Value2 = 0.6696077852558425
$ws.Cells.Item(7,17).Value2 = 121.9229955322309
$ws.Cells.Item(7,18).Value2 = 1097.306959790078
$ws.Cells.Item(7,19).Value2 = 0.1586193581375814
$ws.Cells.Item(7,20).Value2 = 0.1586193581375814

# Row 8
$ws.Cells.Item(8,1).Value2 = "sCs"
$ws.Cells.Item(8,2).Value2 = "Vtn"
$ws.Cells.Item(8,3).Value2 = "Itgb8"
$ws.Cells.Item(8,4).Value2 = "ECs"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 49.26691733333334
$ws.Cells.Item(8,8).Value2 = 147.800752
$ws.Cells.Item(8,9).Value2 = 0.6665842265227129
$ws.Cells.Item(8,10).Value2 = 0.666584226522713
$ws.Cells.Item(8,11).Value2 = 1
$ws.Cells.Item(8,12).Value2 = 0.3333333333333333
$ws.Cells.Item(8,13).Value2 = 0.1126243333333333
$ws.Cells.Item(8,14).Value2 = 0.337873
$ws.Cells.Item(8,15).Value2 = 0.01082936903163217
$ws.Cells.Item(8,16).Value2 = 0.01082936903163217
$ws.Cells.Item(8,17).Value2 = 5.548653720055111
$ws.Cells.Item(8,18).Value2 = 49.93788348049601
$ws.Cells.Item(8,19).Value2 = 0.007218686579679548
$ws.Cells.Item(8,20).Value2 = 0.007218686579679549

# Row 9
$ws.Cells.Item(9,1).Value2 = "sCs"
$ws.Cells.Item(9,2).Value2 = "Vtn"
$ws.Cells.Item(9,3).Value2 = "Itgb8"
$ws.Cells.Item(9,4).Value2 = "FAPs"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 49.26691733333334
$ws.Cells.Item(9,8).Value2 = 147.800752
$ws.Cells.Item(9,9).Value2 = 0.6665842265227129
$ws.Cells.Item(9,10).Value2 = 0.666584226522713
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 3.323421
$ws.Cells.Item(9,14).Value2 = 9.970263
$ws.Cells.Item(9,15).Value2 = 0.3195628457125252
$ws.Cells.Item(9,16).Value2 = 0.3195628457125252
$ws.Cells.Item(9,17).Value2 = 163.734707670864
$ws.Cells.Item(9,18).Value2 = 1473.612369037776
$ws.Cells.Item(9,19).Value2 = 0.2130155523346807
$ws.Cells.Item(9,20).Value2 = 0.2130155523346807

# Row 10
$ws.Cells.Item(10,1).Value2 = "sCs"
$ws.Cells.Item(10,2).Value2 = "Vtn"
$ws.Cells.Item(10,3).Value2 = "Itgb8"
$ws.Cells.Item(10,4).Value2 = "sCs"
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 49.26691733333334
$ws.Cells.Item(10,8).Value2 = 147.800752
$ws.Cells.Item(10,9).Value2 = 0.6665842265227129
$ws.Cells.Item(10,10).Value2 = 0.666584226522713
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 6.963852666666667
$ws.Cells.Item(10,14).Value2 = 20.891558
$ws.Cells.Item(10,15).Value2 = 0.6696077852558425
$ws.Cells.Item(10,16).Value2 = 0.6696077852558425
$ws.Cells.Item(10,17).Value2 = 343.0875536501796
$ws.Cells.Item(10,18).Value2 = 3087.787982851616
$ws.Cells.Item(10,19).Value2 = 0.4463499876083526
$ws.Cells.Item(10,20).Value2 = 0.4463499876083527
